# Rename existing Sheet1 -> Wuhan
$wb = $excel.ActiveWorkbook
$wsWuhan = $wb.Worksheets.Item("Sheet1")
$wsWuhan.Name = "Wuhan"

# Add new sheet "Hubei" positioned right after "Wuhan"; it becomes the active sheet.
$wsHubei = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wsWuhan)
$wsHubei.Name = "Hubei"

# Header row
$wsHubei.Range("A1").Value = "Cofirmed"
$wsHubei.Range("B1").Value = "Death"
$wsHubei.Range("C1").Value = "Recover"

$data = @(
    @(270, 6, 25),
    @(375, 9, 26),
    @(444, 17, 28),
    @(549, 24, 31),
    @(729, 39, 32),
    @(1052, 52, 42),
    @(1423, 76, 44),
    @(2714, 100, 47),
    @(3554, 125, 80),
    @(4586, 162, 90),
    @(5806, 204, 116),
    @(7153, 249, 166),
    @(9074, 294, 215),
    @(11177, 350, 295),
    @(13522, 414, 396),
    @(16678, 479, 520)
)

$r = 2
foreach ($row in $data) {
    $wsHubei.Cells.Item($r, 1).Value = $row[0]
    $wsHubei.Cells.Item($r, 2).Value = $row[1]
    $wsHubei.Cells.Item($r, 3).Value = $row[2]
    $r++
}

# Ensure Hubei ends up the active/selected tab.
$wsHubei.Activate()
$null = $wsHubei.Range("E9").Select()
